# Adds the Unit 44 vocabulary words (U44_01 .. U44_30) to the Vocabulary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column A (UnitID codes) first - matches original authoring order
$ws.Range('A1292').Value = 'U44_01'
$ws.Range('A1293').Value = 'U44_02'
$ws.Range('A1294').Value = 'U44_03'
$ws.Range('A1295').Value = 'U44_04'
$ws.Range('A1296').Value = 'U44_05'
$ws.Range('A1297').Value = 'U44_06'
$ws.Range('A1298').Value = 'U44_07'
$ws.Range('A1299').Value = 'U44_08'
$ws.Range('A1300').Value = 'U44_09'
$ws.Range('A1301').Value = 'U44_10'
$ws.Range('A1302').Value = 'U44_11'
$ws.Range('A1303').Value = 'U44_12'
$ws.Range('A1304').Value = 'U44_13'
$ws.Range('A1305').Value = 'U44_14'
$ws.Range('A1306').Value = 'U44_15'
$ws.Range('A1307').Value = 'U44_16'
$ws.Range('A1308').Value = 'U44_17'
$ws.Range('A1309').Value = 'U44_18'
$ws.Range('A1310').Value = 'U44_19'
$ws.Range('A1311').Value = 'U44_20'
$ws.Range('A1312').Value = 'U44_21'
$ws.Range('A1313').Value = 'U44_22'
$ws.Range('A1314').Value = 'U44_23'
$ws.Range('A1315').Value = 'U44_24'
$ws.Range('A1316').Value = 'U44_25'
$ws.Range('A1317').Value = 'U44_26'
$ws.Range('A1318').Value = 'U44_27'
$ws.Range('A1319').Value = 'U44_28'
$ws.Range('A1320').Value = 'U44_29'
$ws.Range('A1321').Value = 'U44_30'

# Fill remaining columns row by row
# Row 1292 (U44_01)
$ws.Range('B1292').Value = 44
$ws.Range('C1292').Value = 'Thuộc kinh tế'
$ws.Range('D1292').Value = 'Economic'
$ws.Range('E1292').Value = 'Economic growth is projected to be high (dự báo tăng trưởng kinh tế cao)'
$ws.Range('F1292').Value = 'economic growth'
$ws.Range('G1292').Value = 'Adj'

# Row 1293 (U44_02)
$ws.Range('B1293').Value = 44
$ws.Range('C1293').Value = 'Tài chính'
$ws.Range('D1293').Value = 'Financial'
$ws.Range('E1293').Value = 'I need financial support from you'
$ws.Range('F1293').Value = 'Financial support'
$ws.Range('G1293').Value = 'Adj'

# Row 1294 (U44_03)
$ws.Range('B1294').Value = 44
$ws.Range('C1294').Value = 'Tiền mặt'
$ws.Range('D1294').Value = 'Cash'
$ws.Range('E1294').Value = 'Can I pay in cash'
$ws.Range('F1294').Value = 'in cash / bằng tiền mặt'
$ws.Range('G1294').Value = 'N'

# Row 1295 (U44_04)
$ws.Range('B1295').Value = 44
$ws.Range('C1295').Value = 'cổ phiếu'
$ws.Range('D1295').Value = 'Stock'
$ws.Range('E1295').Value = 'You need to research before investing in stock'
$ws.Range('F1295').Value = 'invest in stocks / đầu tư vào cổ phiếu'
$ws.Range('G1295').Value = 'N'

# Row 1296 (U44_05)
$ws.Range('B1296').Value = 44
$ws.Range('C1296').Value = 'Sự đầu tư'
$ws.Range('D1296').Value = 'Investment'
$ws.Range('E1296').Value = 'My parents made a large investment in me'
$ws.Range('F1296').Value = 'make an investment / đầu tư'
$ws.Range('G1296').Value = 'N'

# Row 1297 (U44_06)
$ws.Range('B1297').Value = 44
$ws.Range('C1297').Value = 'Ngành công nghiệp'
$ws.Range('D1297').Value = 'Industry'
$ws.Range('E1297').Value = 'We should network more within the industry'
$ws.Range('F1297').Value = 'within an industry / trong ngành'
$ws.Range('G1297').Value = 'N'

# Row 1298 (U44_07)
$ws.Range('B1298').Value = 44
$ws.Range('C1298').Value = 'Thành lập, thiết lập'
$ws.Range('D1298').Value = 'Set up'
$ws.Range('E1298').Value = 'How do I set up an account on this site?'
$ws.Range('F1298').Value = 'set up an account '
$ws.Range('G1298').Value = 'V'

# Row 1299 (U44_08)
$ws.Range('B1299').Value = 44
$ws.Range('C1299').Value = 'Giao dịch'
$ws.Range('D1299').Value = 'Trade'
$ws.Range('E1299').Value = 'They refuse (từ chối) to trade stock with our company'
$ws.Range('F1299').Value = 'trade stock / giao dịch chứng khoán'
$ws.Range('G1299').Value = 'V'

# Row 1300 (U44_09)
$ws.Range('B1300').Value = 44
$ws.Range('C1300').Value = 'Kiếm tiền'
$ws.Range('D1300').Value = 'Earn'
$ws.Range('E1300').Value = 'He earns a living as a dentist'
$ws.Range('F1300').Value = 'Earn a living / kiếm sống'
$ws.Range('G1300').Value = 'V'

# Row 1301 (U44_10)
$ws.Range('B1301').Value = 44
$ws.Range('C1301').Value = 'Lợi nhuận'
$ws.Range('D1301').Value = 'Profit'
$ws.Range('E1301').Value = 'The business made a huge profit last year'
$ws.Range('F1301').Value = 'make a profit / tạo ra lợi nhuận'
$ws.Range('G1301').Value = 'N'

# Row 1302 (U44_11)
$ws.Range('B1302').Value = 44
$ws.Range('C1302').Value = 'Nhà cung cấp'
$ws.Range('D1302').Value = 'Supplier'
$ws.Range('E1302').Value = 'Japan is a leading supplier of salmon'
$ws.Range('F1302').Value = 'a leading supplier / nhà cung cấp hàng đầu'
$ws.Range('G1302').Value = 'N'

# Row 1303 (U44_12)
$ws.Range('B1303').Value = 44
$ws.Range('C1303').Value = 'Đạt được, có được'
$ws.Range('D1303').Value = 'Obtain'
$ws.Range('E1303').Value = 'After paying the fee, I will obtain a new tourist visa.'
$ws.Range('F1303').Value = 'obtain something / có được một cái gì đó'
$ws.Range('G1303').Value = 'V'

# Row 1304 (U44_13)
$ws.Range('B1304').Value = 44
$ws.Range('C1304').Value = 'Dự định, mục tiêu'
$ws.Range('D1304').Value = 'Aim'
$ws.Range('E1304').Value = 'I exercise a lot with the aim to lose weight.'
$ws.Range('F1304').Value = 'with the aim / với mục tiêu'
$ws.Range('G1304').Value = 'N'

# Row 1305 (U44_14)
$ws.Range('B1305').Value = 44
$ws.Range('C1305').Value = 'Phần trăm'
$ws.Range('D1305').Value = 'Percentage'
$ws.Range('E1305').Value = 'The percentage of your grade is 20%'
$ws.Range('F1305').Value = 'percentage of something / tỷ lệ phần trăm của something'
$ws.Range('G1305').Value = 'N'

# Row 1306 (U44_15)
$ws.Range('B1306').Value = 44
$ws.Range('C1306').Value = 'Giấy phép'
$ws.Range('D1306').Value = 'License'
$ws.Range('E1306').Value = 'You need a license to drive cars'
$ws.Range('F1306').Value = 'license to do something '
$ws.Range('G1306').Value = 'N'

# Row 1307 (U44_16)
$ws.Range('B1307').Value = 44
$ws.Range('C1307').Value = 'Sự thiếu hụt'
$ws.Range('D1307').Value = 'Shortage'
$ws.Range('E1307').Value = 'There is no shortage of food during Tet.'
$ws.Range('F1307').Value = 'a shortage of food / tình trạng thiếu lương thực'
$ws.Range('G1307').Value = 'N'

# Row 1308 (U44_17)
$ws.Range('B1308').Value = 44
$ws.Range('C1308').Value = 'Đóng góp'
$ws.Range('D1308').Value = 'Contribute'
$ws.Range('E1308').Value = 'We all contribute towards economic growth.'
$ws.Range('F1308').Value = 'Contribute to/towards something (đóng góp vào một cái gì đó)'
$ws.Range('G1308').Value = 'V'

# Row 1309 (U44_18)
$ws.Range('B1309').Value = 44
$ws.Range('C1309').Value = 'Thử nghiệm'
$ws.Range('D1309').Value = 'Experiment'
$ws.Range('E1309').Value = 'Scientists conduct experiments regularly'
$ws.Range('F1309').Value = 'conduct an experiment / tiến hành thí nghiệm'
$ws.Range('G1309').Value = 'N'

# Row 1310 (U44_19)
$ws.Range('B1310').Value = 44
$ws.Range('C1310').Value = 'Trận chiến, trận đấu'
$ws.Range('D1310').Value = 'Battle'
$ws.Range('E1310').Value = 'The battle between the fighters(võ sĩ) is on tonight'
$ws.Range('F1310').Value = 'Battle between A and B / trận đấu giữa A và B'
$ws.Range('G1310').Value = 'N'

# Row 1311 (U44_20)
$ws.Range('B1311').Value = 44
$ws.Range('C1311').Value = 'Điều khoản'
$ws.Range('D1311').Value = 'Terms'
$ws.Range('E1311').Value = 'Read the terms and conditions carefully'
$ws.Range('F1311').Value = 'terms and conditions / các điều khoản và điều kiện'
$ws.Range('G1311').Value = 'N'

# Row 1312 (U44_21)
$ws.Range('B1312').Value = 44
$ws.Range('C1312').Value = 'Hàng hóa'
$ws.Range('D1312').Value = 'Goods'
$ws.Range('E1312').Value = 'The store sells luxury goods '
$ws.Range('F1312').Value = 'Luxury goods / hàng hóa cao cấp'
$ws.Range('G1312').Value = 'N'

# Row 1313 (U44_22)
$ws.Range('B1313').Value = 44
$ws.Range('C1313').Value = 'Hiếm'
$ws.Range('D1313').Value = 'Scarce'
$ws.Range('E1313').Value = 'Scarce resources prevent us from continuing'
$ws.Range('F1313').Value = 'scarce resources / Tài nguyên khan hiếm'
$ws.Range('G1313').Value = 'Adj'

# Row 1314 (U44_23)
$ws.Range('B1314').Value = 44
$ws.Range('C1314').Value = 'Mở rộng'
$ws.Range('D1314').Value = 'Extend'
$ws.Range('E1314').Value = 'Can you extend the deadline?'
$ws.Range('F1314').Value = 'Extend the deadline / gia hạn hạn chót'
$ws.Range('G1314').Value = 'V'

# Row 1315 (U44_24)
$ws.Range('B1315').Value = 44
$ws.Range('C1315').Value = 'Khu vực'
$ws.Range('D1315').Value = 'Region'
$ws.Range('E1315').Value = 'Here is the latest news from the Southeast Asian region'
$ws.Range('F1315').Value = 'from the region / từ khu vực'
$ws.Range('G1315').Value = 'N'

# Row 1316 (U44_25)
$ws.Range('B1316').Value = 44
$ws.Range('C1316').Value = 'Lãnh đạo'
$ws.Range('D1316').Value = 'Leader'
$ws.Range('F1316').Value = 'A world leader / nhà lãnh đạo thế giới'
$ws.Range('E1316').Value = 'World leaders will discuss this in the upcoming conference. / Các nhà lãnh đạo thế giới sẽ thảo luận vấn đề này trong hội nghị sắp tới'
$ws.Range('G1316').Value = 'N'

# Row 1317 (U44_26)
$ws.Range('B1317').Value = 44
$ws.Range('C1317').Value = 'Trả hết nợ'
$ws.Range('D1317').Value = 'Pay off'
$ws.Range('E1317').Value = 'Can you pay off the debt?'
$ws.Range('F1317').Value = 'pay of the debt / trả hết nợ'
$ws.Range('G1317').Value = 'V'

# Row 1318 (U44_27)
$ws.Range('B1318').Value = 44
$ws.Range('D1318').Value = 'Scale'
$ws.Range('C1318').Value = 'Quy mô'
$ws.Range('E1318').Value = 'We are distributing(đang phân phối) on a large scale'
$ws.Range('F1318').Value = 'on a large scale / trên một quy mô lớn'
$ws.Range('G1318').Value = 'N'

# Row 1319 (U44_28)
$ws.Range('B1319').Value = 44
$ws.Range('C1319').Value = 'Tiền tệ'
$ws.Range('D1319').Value = 'Currency'
$ws.Range('E1319').Value = 'In VN, US dollar is a foreign currency.'
$ws.Range('F1319').Value = 'foreign currency / ngoại tệ'
$ws.Range('G1319').Value = 'N'

# Row 1320 (U44_29)
$ws.Range('B1320').Value = 44
$ws.Range('C1320').Value = 'Thuế'
$ws.Range('D1320').Value = 'Tax'
$ws.Range('E1320').Value = 'There is a tax on buying products'
$ws.Range('F1320').Value = 'tax on something / thuế cho một thứ gì đó'
$ws.Range('G1320').Value = 'N'

# Row 1321 (U44_30)
$ws.Range('B1321').Value = 44
$ws.Range('C1321').Value = 'Nhấn mạnh'
$ws.Range('D1321').Value = 'Emphasize'
$ws.Range('F1321').Value = 'Emphasize something / nhấn mạnh một cái gì đó'
$ws.Range('E1321').Value = 'The writer emphasized the bravery (quả cảm) of his characters'
$ws.Range('G1321').Value = 'V'

# Move the active selection to the row after the newly added data, matching the saved file state.
[void]$ws.Range('A1322').Select()
